# Mean Values and Box Plot Section
# Rename the sheet from the default "Sheet1" to "DiCE", and move the
# active selection to C27 (where the new Mean Values / Box Plot section
# starts), matching the workbook.xml / sheet1.xml changes in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xl/workbook.xml: <sheet name="Sheet1" .../> -> <sheet name="DiCE" .../>
$ws.Name = "DiCE"

# xl/worksheets/sheet1.xml: <selection activeCell="I13" sqref="I13"/>
#                        -> <selection activeCell="C27" sqref="C27"/>
$ws.Range("C27").Select()
